# Generate Report for Handoff
#
# Updates the "fde9fce8-09c6-47d1-a009-27ee4d22289d.md" row (row 3) on the
# Overview / zh-cn / de-de sheets to reflect that the file is now ready for
# handoff (it was previously "In Translation").

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet (row 3 = fde9fce8...) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Ready for handoff"      # Status
$wsZh.Range("E3").Value = "mt"                     # Priority
$wsZh.Range("H3").Value = "2016-08-20 08:13:50"    # Latest Handoff Datetime

# --- de-de sheet (row 3 = fde9fce8...) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Ready for handoff"      # Status
$wsDe.Range("E3").Value = "mt"                     # Priority
$wsDe.Range("H3").Value = "2016-08-20 08:13:54"    # Latest Handoff Datetime

# --- Overview sheet (row 3 = fde9fce8...) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"        # zh-cn status
$wsOverview.Range("F3").Value = "Ready for handoff"        # de-de status
$wsOverview.Range("G3").Value = "2016-08-20 08:13:54"      # Latest HO Xliff Generate Date (max of the two)

# The longer "Ready for handoff" text no longer fits the previous column
# width, so Excel widened the status columns to fit (autofit). Reproduce
# that widening on the affected columns.
$wsOverview.Range("E1").EntireColumn.ColumnWidth = 16.33   # zh-cn column
$wsOverview.Range("F1").EntireColumn.ColumnWidth = 16.33   # de-de column
$wsZh.Range("C1").EntireColumn.ColumnWidth = 16.33         # Status column
$wsDe.Range("C1").EntireColumn.ColumnWidth = 16.33         # Status column
